$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: human readable column headers (translated from technical identifiers)
$ws.Range("A1").Value = "Número Empresas"
$ws.Range("B1").Value = "Aragón'"
$ws.Range("C1").Value = "Provincia"
$ws.Range("D1").Value = "Mes y año"
$ws.Range("E1").Value = "Dirección provincial (código)"

# Row 2: DSD component identifiers (measures / dimensions)
$ws.Range("A2").Value = "iaest-measure:numero-empresas"
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("C2").Value = "iaest-measure:provincia"
$ws.Range("D2").Value = "iaest-measure:mes-y-ano"
$ws.Range("E2").Value = "null"

# Row 3: component kind (measure / dimension)
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "null"

# Row 4: data type of each component
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "URI-Comunidad"
$ws.Range("C4").Value = "xsd:string"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("E4").Value = "null"
